# This workbook holds weekly price records for "Espinaca" (spinach) at the
# "Terminal Hortofrutícola Agro Chillán" market. The commit re-shuffles the
# per-record fields (Fecha, Calidad, Volumen, Precio mínimo/máximo/promedio,
# Origen, Precio $/Kg) across the existing data rows (rows 2-24), i.e. each
# row receives the values that used to belong to a different row, while the
# constant columns (Mercado ID, Mercado, Región, Codreg, Categoría ID,
# Categoría, Variedad, Unidad de comercialización, Kg o Unidades,
# Clasificación) stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 24

# Columns whose values travel together with each record when rows are
# re-shuffled.
$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

# Maps destination row -> source row (i.e. destRow ends up holding the
# values that sourceRow used to have, before the edit).
$rowMap = @{
    2  = 5
    3  = 15
    4  = 17
    5  = 6
    6  = 20
    7  = 21
    8  = 11
    9  = 22
    10 = 13
    11 = 12
    12 = 23
    13 = 19
    14 = 7
    15 = 24
    16 = 14
    17 = 10
    18 = 8
    19 = 2
    20 = 9
    21 = 16
    22 = 3
    23 = 4
    24 = 18
}

# First snapshot all the current ("before") values for the columns that move,
# so the writes below (which happen in-place) never read already-overwritten
# data. (Uses Value2 for reading - it reliably returns the underlying
# scalar in this host, whereas Value's getter does not.)
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the permuted values into each destination row.
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
